$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.86145555973053
$ws.Range("B1").Value = 1.999409437179565
$ws.Range("C1").Value = 2.159294366836548
$ws.Range("D1").Value = 3.091058254241943
$ws.Range("E1").Value = 2.878909349441528
